# Fix contact information missing from short resumes
# Insert a centered contact-info paragraph right after the name heading
# ("Dheeraj Chand") and before the "PROFESSIONAL SUMMARY" section.

$d = $word.ActiveDocument

# The name/title paragraph is the document's first paragraph.
$titlePara = $d.Paragraphs.Item(1)

# Create a brand-new (initially empty) paragraph right after it. This
# paragraph does not yet have the text we want, and it temporarily
# inherits the title run's formatting (bold / 14pt) -- that gets
# replaced wholesale below via InsertXML so the final paragraph has a
# clean, non-bold, default-size run just like the rest of the body text.
$titlePara.Range.InsertParagraphAfter()
$contactPara = $d.Paragraphs.Item(2)
$contactRange = $contactPara.Range

$contactXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$contactRange.InsertXML($contactXml)
